$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-32 down to 10-33
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new record
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = "Macroferia Regional de Talca"
$ws.Range("C9").Value = "Maule"
$ws.Range("D9").Value = 44481
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100107
$ws.Range("H9").Value = "Otros"
$ws.Range("I9").Value = 100107002
$ws.Range("J9").Value = "Chirimoya"
$ws.Range("K9").Value = "Cultivar IV Región"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 150
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 25000
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("R9").Value = "Provincia de Limarí"
$ws.Range("S9").Value = 2500
$ws.Range("T9").Value = 10
